$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.207.93'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.477.87'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').Value = "'577.98"
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('D6').Value = "'147.04"
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '2.477.78'
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('E14').Value = '  +4.84%  '
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '2.928.30'
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('D17').Value = '63.222.21'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').Value = '2.483.82'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').Value = "'8.24"
$ws.Range('E19').Value = '  +3.67%  '
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').Value = "'329.70"
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').Value = "'2.27"
$ws.Range('E22').Value = '  +10.12%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = "'0.999"
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  +1.14%  '
$ws.Range('D26').Value = "'674.61"
$ws.Range('E26').Value = '  +6.00%  '
$ws.Range('D27').Value = "'9.73"
$ws.Range('E27').Value = '  +13.82%  '
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('D29').Value = '2.634.49'
$ws.Range('E29').Value = '  +3.79%  '
$ws.Range('E30').Value = '  -9.58%  '
$ws.Range('D32').Value = "'8.07"
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('D34').Value = "'0.134"
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('E35').Value = '  +4.03%  '
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('E39').Value = '  -0.67%  '
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('D41').Value = "'151.10"
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('D42').Value = "'2.75"
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '0.0₆0311'
$ws.Range('E45').Value = '  +9.25%  '
$ws.Range('D46').Value = "'154.18"
$ws.Range('E46').Value = '  +6.12%  '
$ws.Range('E47').Value = '  +18.67%  '
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').Value = "'0.0513"
$ws.Range('E51').Value = '  -0.78%  '
